$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

Set-CellText $ws.Range('D2') '60.919.10'
Set-CellText $ws.Range('E2') '  -4.80%  '
Set-CellText $ws.Range('D3') '2.961.43'
Set-CellText $ws.Range('E3') '  -4.38%  '
Set-CellText $ws.Range('E4') '  -0.04%  '
Set-CellText $ws.Range('D5') '541.04'
Set-CellText $ws.Range('E5') '  -0.66%  '
Set-CellText $ws.Range('D6') '130.34'
Set-CellText $ws.Range('E6') '  -7.19%  '
Set-CellText $ws.Range('D8') '2.957.16'
Set-CellText $ws.Range('E8') '  -4.37%  '
Set-CellText $ws.Range('E9') '  -1.94%  '
Set-CellText $ws.Range('E10') '  -7.99%  '
Set-CellText $ws.Range('D11') '5.83'
Set-CellText $ws.Range('E11') '  -10.77%  '
Set-CellText $ws.Range('D12') '0.439'
Set-CellText $ws.Range('E12') '  -4.15%  '
Set-CellText $ws.Range('D13') '0.0000217'
Set-CellText $ws.Range('E13') '  -4.31%  '
Set-CellText $ws.Range('D14') '33.45'
Set-CellText $ws.Range('E14') '  -4.21%  '
Set-CellText $ws.Range('D15') '3.445.08'
Set-CellText $ws.Range('E15') '  -4.07%  '
Set-CellText $ws.Range('B16') 'TRON'
Set-CellText $ws.Range('C16') 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
Set-CellText $ws.Range('D16') '0.109'
Set-CellText $ws.Range('E16') '  -3.36%  '
Set-CellText $ws.Range('B17') 'WrappedBTC'
Set-CellText $ws.Range('C17') 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-CellText $ws.Range('D17') '60.939.67'
Set-CellText $ws.Range('E17') '  -4.89%  '
Set-CellText $ws.Range('D18') '2.968.55'
Set-CellText $ws.Range('E18') '  -4.02%  '
Set-CellText $ws.Range('E19') '  -2.51%  '
Set-CellText $ws.Range('D20') '461.00'
Set-CellText $ws.Range('E20') '  -4.06%  '
Set-CellText $ws.Range('D21') '12.96'
Set-CellText $ws.Range('E21') '  -3.78%  '
Set-CellText $ws.Range('D22') '0.658'
Set-CellText $ws.Range('E22') '  -6.21%  '
Set-CellText $ws.Range('D23') '6.86'
Set-CellText $ws.Range('E23') '  -3.66%  '
Set-CellText $ws.Range('D24') '79.07'
Set-CellText $ws.Range('E24') '  +0.04%  '
Set-CellText $ws.Range('D25') '11.80'
Set-CellText $ws.Range('D26') '0.998'
Set-CellText $ws.Range('E26') '  -0.26%  '
Set-CellText $ws.Range('D27') '2.66'
Set-CellText $ws.Range('E27') '  -2.38%  '
Set-CellText $ws.Range('D28') '7.50'
Set-CellText $ws.Range('E28') '  -7.30%  '
Set-CellText $ws.Range('E29') '  +0.18%  '
Set-CellText $ws.Range('E30') '  -2.68%  '
Set-CellText $ws.Range('D31') '25.06'
Set-CellText $ws.Range('E31') '  -5.03%  '
Set-CellText $ws.Range('E32') '  -4.80%  '
Set-CellText $ws.Range('E33') '  -4.98%  '
Set-CellText $ws.Range('E34') '  -1.40%  '
Set-CellText $ws.Range('D35') '54.12'
Set-CellText $ws.Range('E35') '  -5.54%  '
Set-CellText $ws.Range('D36') '5.77'
Set-CellText $ws.Range('E36') '  -4.36%  '
Set-CellText $ws.Range('D37') '441.10'
Set-CellText $ws.Range('E37') '  -11.51%  '
Set-CellText $ws.Range('D38') '3.126.02'
Set-CellText $ws.Range('E38') '  -3.96%  '
Set-CellText $ws.Range('D39') '0.0778'
Set-CellText $ws.Range('E39') '  -3.04%  '
Set-CellText $ws.Range('B40') 'Kaspa'
Set-CellText $ws.Range('C40') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-CellText $ws.Range('D40') '0.116'
Set-CellText $ws.Range('E40') '  -3.54%  '
Set-CellText $ws.Range('B41') 'VeChain'
Set-CellText $ws.Range('C41') 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-CellText $ws.Range('D41') '0.0372'
Set-CellText $ws.Range('E41') '  -8.25%  '
Set-CellText $ws.Range('D42') '7.97'
Set-CellText $ws.Range('E42') '  -2.06%  '
Set-CellText $ws.Range('E44') '  -14.22%  '
Set-CellText $ws.Range('D45') '25.29'
Set-CellText $ws.Range('E45') '  -0.11%  '
Set-CellText $ws.Range('E46') '  -6.94%  '
Set-CellText $ws.Range('E47') '  -2.80%  '
Set-CellText $ws.Range('D48') '115.56'
Set-CellText $ws.Range('E48') '  -7.14%  '
Set-CellText $ws.Range('D49') '1.91'
Set-CellText $ws.Range('E49') '  -6.78%  '
Set-CellText $ws.Range('E50') '  +7.04%  '
Set-CellText $ws.Range('D51') '0.0₃0474'
Set-CellText $ws.Range('E51') '  -11.41%  '
